$d = $word.ActiveDocument

# 1) Fix the email address in the header/contact line
$found = $d.Content.Find.Execute("joe@wanat.com", $true, $false, $false, $false, $false,
                                  $true, 1, $false, "joe@bla.com", 2)

# 2) Remove the hanging/negative first-line indent ("w:first-line=-720") that was
#    set on every bulleted ("<tab>\u2022<tab>...") paragraph, leaving just the
#    720-twip left indent behind.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("•")) {
        $p.Format.FirstLineIndent = $null
    }
}
